$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -1
$ws.Range("B1").Value = 3.769417524337769
$ws.Range("C1").Value = 2.079323053359985
$ws.Range("D1").Value = 1.630290031433105
$ws.Range("E1").Value = 1.48291015625
